# Auto-generated Excel COM-interop script to update cryptos price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values, forcing text storage so values like
# "1.00" or "25.02" are not silently coerced into numbers by Excel.
$priceUpdates = @{
    2 = '64.456.39'
    3 = '3.505.96'
    5 = '587.22'
    6 = '135.80'
    7 = '3.506.64'
    9 = '0.488'
    10 = '0.125'
    12 = '0.377'
    13 = '4.099.53'
    14 = '0.0000182'
    16 = '3.501.98'
    17 = '64.439.64'
    18 = '25.02'
    19 = '10.05'
    20 = '5.65'
    21 = '13.81'
    22 = '385.47'
    24 = '3.643.80'
    25 = '74.16'
    26 = '1.00'
    30 = '7.51'
    31 = '1.00'
    34 = '3.523.39'
    37 = '23.60'
    38 = '5.32'
    39 = '1.56'
    40 = '6.86'
    41 = '163.58'
    42 = '0.0785'
    43 = '0.809'
    44 = '26.06'
    46 = '41.90'
    47 = '1.22'
    48 = '4.43'
    50 = '2.476.82'
    51 = '6.79'
}
foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}

# Update Volume(1h) (column E) percentage-change text values.
$volumeUpdates = @{
    2 = '  -0.26%  '
    3 = '  -0.09%  '
    4 = '  -0.04%  '
    5 = '  +0.28%  '
    6 = '  +2.77%  '
    7 = '  -0.10%  '
    8 = '  +0.00%  '
    9 = '  -0.36%  '
    10 = '  +0.58%  '
    11 = '  -0.83%  '
    12 = '  -2.58%  '
    13 = '  -0.24%  '
    14 = '  +1.21%  '
    15 = '  +1.35%  '
    16 = '  -0.34%  '
    17 = '  -0.27%  '
    18 = '  -9.50%  '
    19 = '  +1.33%  '
    20 = '  -0.71%  '
    21 = '  -2.62%  '
    22 = '  -1.29%  '
    23 = '  -1.28%  '
    24 = '  -0.28%  '
    25 = '  +0.24%  '
    26 = '  +0.02%  '
    27 = '  +1.67%  '
    28 = '  +3.69%  '
    29 = '  -0.17%  '
    30 = '  +1.02%  '
    31 = '  +0.10%  '
    32 = '  -1.06%  '
    33 = '  +0.54%  '
    34 = '  +0.20%  '
    35 = '  -0.01%  '
    36 = '  +0.98%  '
    37 = '  -1.69%  '
    38 = '  +0.94%  '
    39 = '  -2.13%  '
    40 = '  -1.37%  '
    41 = '  -4.35%  '
    42 = '  -2.59%  '
    43 = '  -0.55%  '
    44 = '  -1.50%  '
    45 = '  -0.04%  '
    46 = '  -0.49%  '
    47 = '  +0.18%  '
    48 = '  +0.84%  '
    49 = '  +0.07%  '
    50 = '  +0.68%  '
    51 = '  -1.38%  '
}
foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}
